$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.000680088996887
$ws.Range("B1").Value = 1.799864172935486
$ws.Range("C1").Value = 5.05289888381958
$ws.Range("D1").Value = 1.428756713867188
$ws.Range("E1").Value = 1.344442129135132
